$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 91 (A column) was the previous "last row" and used a date-only format.
# Now that a new row is appended below it, it becomes a regular data row,
# so switch it to the same datetime format used by the other interior rows.
$ws.Cells.Item(91, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data as row 92.
$ws.Cells.Item(92, 1).Value = 45831
$ws.Cells.Item(92, 2).Value = 387
$ws.Cells.Item(92, 3).Value = 392
$ws.Cells.Item(92, 4).Value = 395

# Row 92 becomes the new "last row" and gets the date-only format that row 91 used to have.
$ws.Cells.Item(92, 1).NumberFormat = "YYYY-MM-DD"
